# Trade #101 closed at 2026-02-17 15:57:59 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.95   # Current Capital
$wsSummary.Range("B4").Value = -0.06     # Total P&L $
$wsSummary.Range("B5").Value = -0.01     # Total P&L %
$wsSummary.Range("B6").Value = 101       # Total Trades
$wsSummary.Range("B7").Value = 37        # Winning Trades
$wsSummary.Range("B9").Value = 36.63     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.95
$wsStatus.Range("D4").Value = 101
$wsStatus.Range("E4").Value = -0.06
$wsStatus.Range("F4").Value = -0.05
$wsStatus.Range("G4").Value = 36.63

# ---------------------------------------------------------------------
# Append trade #101 as a new row (row 102) to both the "All Trades" and
# "MarketMaking" logs - they mirror each other.
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A102").Value = 101

    # Write the date as literal text instead of letting it be parsed into
    # a date serial: force a text number format before assigning, then
    # drop the formatting override so the cell keeps the default style.
    $ws.Range("B102").NumberFormat = "@"
    $ws.Range("B102").Value = "2026-02-17"
    $ws.Range("B102").ClearFormats()

    $ws.Range("C102").Value = "15:57:52"
    $ws.Range("D102").Value = "MarketMaking"
    $ws.Range("E102").Value = "DOWN"
    $ws.Range("F102").Value = 0.86
    $ws.Range("G102").Value = 0.92
    $ws.Range("H102").Value = "CLOSED"
    $ws.Range("I102").Value = 6.9767
    $ws.Range("J102").Value = 0.06
    $ws.Range("K102").Value = 99.95
    $ws.Range("L102").Value = 0
    $ws.Range("M102").Value = 0
    $ws.Range("N102").Value = 0.6
    $ws.Range("O102").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P102").Value = "early_exit"
    $ws.Range("Q102").Value = 0.14
}
